# Regenerate save_data: column G ("K" - strikeouts) values updated to use K
# instead of the old Strike# derived figures. Only column G values change;
# all other columns/rows are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"  = 1
    "G4"  = 0
    "G5"  = 1
    "G6"  = 1
    "G7"  = 2
    "G8"  = 1
    "G9"  = 1
    "G10" = 1
    "G11" = 1
    "G12" = 1
    "G13" = 0
    "G14" = 2
    "G15" = 0
    "G18" = 2
    "G19" = 1
    "G20" = 1
    "G21" = 0
    "G22" = 0
    "G23" = 0
    "G24" = 0
    "G25" = 1
    "G26" = 0
    "G27" = 1
    "G28" = 1
    "G29" = 3
    "G30" = 2
    "G31" = 1
    "G32" = 1
    "G33" = 2
    "G34" = 1
    "G35" = 1
    "G36" = 0
    "G37" = 1
    "G38" = 1
    "G39" = 2
    "G40" = 0
    "G41" = 0
    "G42" = 0
    "G43" = 2
    "G44" = 1
    "G45" = 1
    "G46" = 0
    "G48" = 2
    "G49" = 1
    "G50" = 1
    "G51" = 3
    "G52" = 2
    "G53" = 2
    "G54" = 0
    "G55" = 0
    "G56" = 1
    "G57" = 1
    "G58" = 1
    "G60" = 1
    "G61" = 0
    "G62" = 2
    "G63" = 1
    "G64" = 2
    "G68" = 0
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
